$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): D1 text changes from "Backend Developer" to "Email ID" ---
$ws.Range("D1").Value = "Email ID"

# --- Data rows: replace rows 2-4 with new employee data ---
# (write names first, then ids/emails/roles, then skills - matches original authoring order)
$ws.Range("B2").Value = "Ashu"
$ws.Range("C2").Value = "Das"
$ws.Range("B3").Value = "Nilesh"
$ws.Range("C3").Value = "Ghosh"
$ws.Range("B4").Value = "Vijay"
$ws.Range("C4").Value = "Tripathi"

$ws.Range("A2").Value = 458789
$ws.Range("D2").Value = "asdasdasd@temp.com"
$ws.Range("E2").Value = "Fullstack Developer"

$ws.Range("A3").Value = 589698
$ws.Range("D3").Value = "dgdfgdfg@temp.com"
$ws.Range("E3").Value = "Fullstack Developer"

$ws.Range("A4").Value = 215468
$ws.Range("D4").Value = "tertertert@temp.com"
$ws.Range("E4").Value = "Fullstack Developer"

$ws.Range("F2").Value = ".Net, SQL, Angular"
$ws.Range("F3").Value = ".Net, SQL, Angular"
$ws.Range("F4").Value = "Angular, React"

# --- New font formatting for the "Project Role" column (E2:E4) ---
$font = $ws.Range("E2").Font
$font.Color = 7901646  # BGR long value equal to RGB(0xCE,0x91,0x78) == ARGB FFCE9178
$font.Size = 12
$font.Name = "Consolas"
$ws.Range("E2").VerticalAlignment = -4108  # xlCenter

$ws.Range("E2").Copy()
$ws.Range("E3:E4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row heights for data rows 2-4 ---
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75
$ws.Rows.Item(4).RowHeight = 15.75

# --- New trailing empty row (row 5) ---
$ws.Range("D5").Style = "Hyperlink"
$ws.Range("E5").WrapText = $true
$ws.Range("E5").VerticalAlignment = -4108  # xlCenter
$ws.Range("F5").WrapText = $true
$ws.Range("F5").VerticalAlignment = -4108  # xlCenter

# --- Column F width ---
$ws.Columns.Item(6).ColumnWidth = 31.5

# --- Selection change ---
$ws.Range("E9").Select()
